$d = $word.ActiveDocument

function Merge-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph: Title ---
# "a" + " short and descriptive name for the bug" + " which" + " " -> "a short and descriptive name for the bug which "
Merge-Text "a short and descriptive name for the bug which " "a short and descriptive name for the bug which "
# "summari" + "s" + "e" + "s" + " the issue" -> "summarises the issue"
Merge-Text "summarises the issue" "summarises the issue"

# --- Paragraph: Bug Report ID ---
# "- a " + "unique identifier assigned to " + "each" + " bug report" -> "- a unique identifier assigned to each bug report"
Merge-Text "- a unique identifier assigned to each bug report" "- a unique identifier assigned to each bug report"

# --- Paragraph: Linked Test Case ID ---
# Text change: TC_X_00x -> TC_X_0xx, plus merge of the trailing description runs
Merge-Text "Linked Test Case ID: TC_X_00x" "Linked Test Case ID: TC_X_0xx"
Merge-Text "- the ID of the test case associated with the bug (if applicable)" "- the ID of the test case associated with the bug (if applicable)"

# --- Paragraph: Preconditions ---
# "any" stays separate; " " + "necessary configurations" -> " necessary configurations"
Merge-Text " necessary configurations" " necessary configurations"

# --- Paragraph: Steps to Reproduce (intro) ---
# " for " + "replicating the bug" -> " for replicating the bug"
Merge-Text " for replicating the bug" " for replicating the bug"

# --- Paragraph: the second step ---
Merge-Text "the second step " "the second step "

# --- Paragraph: the third step ---
Merge-Text "the third step" "the third step"

# --- Paragraph: … (steps) ---
Merge-Text "… " "… "

# --- Paragraph: Test Data ---
Merge-Text "any specific input data or credentials required for testing" "any specific input data or credentials required for testing"

# --- Paragraph: Expected Result ---
Merge-Text "the intended behaviour of the application if it is working as expected" "the intended behaviour of the application if it is working as expected"

# --- Paragraph: the second expected outcome ---
Merge-Text "the second expected outcome" "the second expected outcome"

# --- Paragraph: the third expected outcome ---
Merge-Text "the third expected outcome" "the third expected outcome"

# --- Paragraph: Reproducibility ---
# ", conditional " -> ", conditional) " (genuine text insertion of a closing paren)
Merge-Text ", conditional " ", conditional) "
